$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 39
    3  = 46
    4  = 32
    5  = 27
    6  = 125
    7  = 83
    8  = 28
    9  = 40
    10 = 88
    11 = 175
    12 = 24
    13 = 56
    14 = 91
    15 = 7
    16 = 231
    17 = 1
    18 = 115
    19 = 159
    20 = 287
    21 = 47
    23 = 23
    24 = 58
    25 = 43
    26 = 57
    27 = 98
    28 = 70
    29 = 85
    30 = 111
    31 = 19
    32 = 31
    33 = 99
    34 = 60
    35 = 110
    36 = 66
    37 = 53
    38 = 86
    39 = 138
    40 = 55
    41 = 87
    42 = 137
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
